# fix : dictionary error
# Update bird data table: simplify Korean names, replace habitat/population
# columns with endangered-status descriptions, refresh headers, resize columns,
# and move the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range('D1').Value = '멸종 위기 등급'
$ws.Range('F1').Value = '포인트'
$ws.Range('G1').Value = '확률'

# --- Row 2: Crow / 까마귀 ---
$ws.Range('B2').Value = '까마귀'
$ws.Range('C2').Value = '전세계'
$ws.Range('D2').Value = '해당 없음.'
$ws.Range('E2').Value = '성체 기준 몸 길이는 48~52cm이며, 무게는 약 300~600g으로 대게 수컷이 암컷보다 크다. 매우 지능적이며 인간 환경에 적응할 수 있다.'

# --- Row 3: Golden Eagle / 검독수리 ---
$ws.Range('C3').Value = '유라시아, 북아메리카, 일부 아프리카 지역'
$ws.Range('D3').Value = '대한민국 멸종위기 I급.'
$ws.Range('E3').Value = '전체적으로 어두운 갈색을 띄며, 머리와 목에 더 연한 깃털이 있다. 익장이 2m에 달하며 최대 시속 240km로 비행할 수 있다. 먹이는 주로 토끼, 다람쥐, 청설모, 거북, 뱀 등이 있다.'

# --- Row 4: Great Horned Owl / 수리부엉이 ---
$ws.Range('B4').Value = '수리부엉이'
$ws.Range('C4').Value = '유라시아 대부분 지역, 유럽 전 지역'
$ws.Range('D4').Value = '해당 없음.'
$ws.Range('E4').Value = '키는 위아래로 60 ~ 75cm, 익장 131 ~ 188cm. 매우 튼튼한 다리, 날카로운 발톱을 가지고 있다. 몸에 있는 부드러운 솜털이 소리를 흡수하기 때문에 거의 무소음에 가까운 비행이 가능하다. '

# --- Row 5: Pigeon / 집비둘기 ---
$ws.Range('C5').Value = '전세계'
$ws.Range('D5').Value = '해당 없음.'
$ws.Range('E5').Value = '1년 내내 번식이 가능하다. 최고 시속 112km를 자랑하며 기억력이 높고, 인간의 얼굴을 구별할 수 있다.'

# --- Row 6: Seagull / 갈매기 ---
$ws.Range('B6').Value = '갈매기'
$ws.Range('C6').Value = '바닷가'
$ws.Range('D6').Value = '관심 필요.'
$ws.Range('E6').Value = '바닷가에 주로 서식하며 그 외에도 습지, 내륙의 호수 등에서도 흔히 볼 수 있다. 몸길이는 약 40cm이고 날개를 폈을 때는 120cm 정도이다.'

# --- Row 7: Sparrow / 참새 ---
$ws.Range('B7').Value = '참새'
$ws.Range('C7').Value = '전세계에 걸쳐 분포. 종마다 상이'
$ws.Range('D7').Value = '해당 없음.'
$ws.Range('E7').Value = '다양한 계통으로 종이 나뉘어져 있으며 잡식성으로 곡식 낟알이나 벌레를 잡아먹는다.'

# --- Column widths (C/D/E widened for the longer descriptive text) ---
$ws.Columns.Item(3).ColumnWidth = 47.857142857142854
$ws.Columns.Item(4).ColumnWidth = 40.42857142857143
$ws.Columns.Item(5).ColumnWidth = 78

# --- Selection moves to D8 ---
$ws.Range('D8').Select()
